# Progress-report update: append a new "Xây dựng mô hình ... " status block
# right after the "Hoàn thành cơ bản công việc thứ nhất (Nghiên cứu lý
# thuyết)." bullet, and before the trailing indented/tab paragraph.

$d = $word.ActiveDocument

$anchorText = "Hoàn thành cơ bản công việc thứ nhất (Nghiên cứu lý thuyết)."

$findRange = $d.Content
$found = $findRange.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Anchor paragraph not found"
}

# Build a *fresh* zero-length Range object positioned one character before
# the end of the match (i.e. strictly inside the trailing "." run, not
# sitting exactly on the paragraph mark). Collapsing the Find range itself
# (or inserting exactly at a paragraph-mark boundary) causes InsertXML to
# swallow the neighbouring paragraph, so we avoid that by re-creating the
# Range from plain integer offsets instead.
$insertAt = $findRange.End - 1
$ins = $d.Range($insertAt, $insertAt)

$xmlPayload = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="vi-VN"/></w:rPr><w:t>Xây dựng mô hình</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> từ</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="vi-VN"/></w:rPr><w:t>16</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>/11 đế</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">n </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="vi-VN"/></w:rPr><w:t>29</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>/11</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:bCs/><w:lang w:val="vi-VN"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:bCs/><w:lang w:val="vi-VN"/></w:rPr><w:t>Đã thực hiện: Xây dựng sơ đồ lớp, lưu đồ xử lý, viết báo cáo cho chương 3 (Hiện thực hóa nghiên cứu).</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$ins.InsertXML($xmlPayload)

Write-Output "inserted at $insertAt"
